$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 189, shifting existing rows 189-204 down to 190-205.
$ws.Rows.Item(189).Insert()

# Populate the new row 189 with the new record.
$ws.Cells.Item(189, 1).Value = 7
$ws.Cells.Item(189, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(189, 3).Value = "Ñuble"
$ws.Cells.Item(189, 4).Value = 45223
$ws.Cells.Item(189, 5).Value = 16
$ws.Cells.Item(189, 6).Value = "Fruta"
$ws.Cells.Item(189, 7).Value = 100108
$ws.Cells.Item(189, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(189, 9).Value = 100108002
$ws.Cells.Item(189, 10).Value = "Mango"
$ws.Cells.Item(189, 11).Value = "Sin especificar"
$ws.Cells.Item(189, 12).Value = "Primera"
$ws.Cells.Item(189, 13).Value = 50
$ws.Cells.Item(189, 14).Value = 10000
$ws.Cells.Item(189, 15).Value = 10000
$ws.Cells.Item(189, 16).Value = 10000
$ws.Cells.Item(189, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(189, 18).Value = "Brasil"
$ws.Cells.Item(189, 19).Value = 2500
$ws.Cells.Item(189, 20).Value = 4
